# Costos 20100927 - add Costos/QC/Burndown charts update for the 27/09 delivery.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

# --- Update the raw sprint-1 measurements; B4:B7/ C/D formulas recalc automatically ---
$ws.Range("B3").Value = 72
$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 22.5

# --- Row 2 (table header) goes back to the default row height ---
$ws.Rows(2).AutoFit()

# --- Remove the chart's explicit "Costo" title (now relies on auto title deleted) ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.HasTitle = $false

# --- Keep the chart series bound/refreshed against the (now recalculated) source data ---
$chart.SetSourceData($ws.Range("A2:D7"))

# --- Move the selection cursor like the author left it ---
$ws.Range("C14").Select()
